$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 222.2
$ws.Range("I6").Value = 35.666668
$ws.Range("J6").Value = 502
$ws.Range("K6").Value = 107.000004
$ws.Range("L6").Value = 1506
$ws.Range("M6").Value = 4.999995999999996
$ws.Range("N6").Value = -1730
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("K11").Value = 1
$ws.Range("M11").Value = 139
$ws.Range("H26").Value = 5350
$ws.Range("J26").Value = 5350
$ws.Range("L26").Value = 5350
$ws.Range("N26").Value = -6038
$ws.Range("H40").Value = 8774.375
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 8774.375
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 8774.375
$ws.Range("N40").Value = -9124.375
$ws.Range("H47").Value = 11000
$ws.Range("I47").Value = 14000
$ws.Range("J47").Value = 8000
$ws.Range("K47").Value = 14000
$ws.Range("L47").Value = 8000
$ws.Range("M47").Value = -13028
$ws.Range("N47").Value = -9944
$ws.Range("H51").Value = 4712.25
$ws.Range("I51").Value = 3000
$ws.Range("J51").Value = 5283
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 5283
$ws.Range("M51").Value = -2516
$ws.Range("N51").Value = -6251
$ws.Range("H62").Value = 4899.5
$ws.Range("I62").Value = 4899
$ws.Range("K62").Value = 4899
$ws.Range("M62").Value = -4275
$ws.Range("H64").Value = 7309.6206
$ws.Range("I64").Value = 6071.077
$ws.Range("J64").Value = 8315.9375
$ws.Range("K64").Value = 6071.077
$ws.Range("L64").Value = 8315.9375
$ws.Range("M64").Value = -5823.077
$ws.Range("N64").Value = -8811.9375
$ws.Range("H65").Value = 4899.5
$ws.Range("I65").Value = 4899
$ws.Range("K65").Value = 24495
$ws.Range("M65").Value = -21375
$ws.Range("H67").Value = 7309.6206
$ws.Range("I67").Value = 6071.077
$ws.Range("J67").Value = 8315.9375
$ws.Range("K67").Value = 6071.077
$ws.Range("L67").Value = 8315.9375
$ws.Range("M67").Value = -5213.077
$ws.Range("N67").Value = -10031.9375
$ws.Range("H141").Value = 3997.2144
$ws.Range("I141").Value = 3920.8462
$ws.Range("K141").Value = 11762.5386
$ws.Range("M141").Value = -6582.5386
$ws.Range("M40").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17255698
$ws.Range("I32").Value = 19619214
$ws.Range("K32").Value = 19619214
$ws.Range("M32").Value = -19618927
$ws.Range("H63").Value = 2071.611
$ws.Range("I63").Value = 2115
$ws.Range("J63").Value = 1724.5
$ws.Range("K63").Value = 2115
$ws.Range("L63").Value = 1724.5
$ws.Range("M63").Value = -1429
$ws.Range("N63").Value = -3096.5
$ws.Range("H66").Value = 2071.611
$ws.Range("I66").Value = 2115
$ws.Range("J66").Value = 1724.5
$ws.Range("K66").Value = 10575
$ws.Range("L66").Value = 8622.5
$ws.Range("M66").Value = -7143
$ws.Range("N66").Value = -15486.5
$ws.Range("H110").Value = 3187.4546
$ws.Range("I110").Value = 3187.4546
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 3187.4546
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -1142.4546
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3247.12
$ws.Range("I105").Value = 1748.6666
$ws.Range("J105").Value = 4630.3076
$ws.Range("K105").Value = 1748.6666
$ws.Range("L105").Value = 4630.3076
$ws.Range("M105").Value = -1.666600000000017
$ws.Range("N105").Value = -8124.3076

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1666906.6
$ws.Range("J22").Value = 5000000
$ws.Range("L22").Value = 5000000
$ws.Range("N22").Value = -5000700
$ws.Range("H99").Value = 9763567
$ws.Range("I99").Value = 2038461.6
$ws.Range("J99").Value = 18190956
$ws.Range("K99").Value = 2038461.6
$ws.Range("L99").Value = 18190956
$ws.Range("M99").Value = -2036963.6
$ws.Range("N99").Value = -18193952
$ws.Range("H122").Value = 5205209.5
$ws.Range("I122").Value = 10405745
$ws.Range("J122").Value = 4673.4707
$ws.Range("K122").Value = 31217235
$ws.Range("L122").Value = 14020.4121
$ws.Range("M122").Value = -31214785
$ws.Range("N122").Value = -18920.4121
$ws.Range("H126").Value = 9763567
$ws.Range("I126").Value = 2038461.6
$ws.Range("J126").Value = 18190956
$ws.Range("K126").Value = 6115384.800000001
$ws.Range("L126").Value = 54572868
$ws.Range("M126").Value = -6112914.800000001
$ws.Range("N126").Value = -54577808

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 143.375
$ws.Range("J2").Value = 236.25
$ws.Range("L2").Value = 1417.5
$ws.Range("N2").Value = -1643.5
$ws.Range("H16").Value = 500
$ws.Range("I16").Value = 500
$ws.Range("K16").Value = 1500
$ws.Range("M16").Value = -1327
$ws.Range("H38").Value = 1700.75
$ws.Range("I38").Value = 43.2
$ws.Range("J38").Value = 2884.7144
$ws.Range("K38").Value = 129.6
$ws.Range("L38").Value = 8654.143199999999
$ws.Range("M38").Value = 217.4
$ws.Range("N38").Value = -9348.143199999999
$ws.Range("H68").Value = 940
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 940
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2820
$ws.Range("N68").Value = -4442
$ws.Range("H71").Value = 940
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 940
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 8460
$ws.Range("N71").Value = -16572
$ws.Range("H129").Value = 1246
$ws.Range("I129").Value = 858.75
$ws.Range("J129").Value = 2795
$ws.Range("K129").Value = 2576.25
$ws.Range("L129").Value = 8385
$ws.Range("M129").Value = 2423.75
$ws.Range("N129").Value = -18385
$ws.Range("H131").Value = 5649.846
$ws.Range("I131").Value = 3344.8
$ws.Range("K131").Value = 10034.4
$ws.Range("M131").Value = -4994.400000000001
$ws.Range("H132").Value = 2467
$ws.Range("I132").Value = 900.75
$ws.Range("J132").Value = 5599.5
$ws.Range("K132").Value = 8106.75
$ws.Range("L132").Value = 50395.5
$ws.Range("M132").Value = -5576.75
$ws.Range("N132").Value = -55455.5
$ws.Range("H134").Value = 6213.385
$ws.Range("I134").Value = 2444
$ws.Range("K134").Value = 7332
$ws.Range("M134").Value = -2262
$ws.Range("H137").Value = 2433.8
$ws.Range("I137").Value = 1415.1428
$ws.Range("K137").Value = 4245.428400000001
$ws.Range("M137").Value = 854.5715999999993
$ws.Range("H139").Value = 23260692
$ws.Range("I139").Value = 37039732
$ws.Range("K139").Value = 111119196
$ws.Range("M139").Value = -111114056
$ws.Range("M68").ClearContents()
$ws.Range("M71").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 33224.87
$ws.Range("I70").Value = 62458.844
$ws.Range("K70").Value = 62458.844
$ws.Range("M70").Value = -62188.844
$ws.Range("H73").Value = 33224.87
$ws.Range("I73").Value = 62458.844
$ws.Range("K73").Value = 62458.844
$ws.Range("M73").Value = -61522.844
$ws.Range("H80").Value = 12866259
$ws.Range("I80").Value = 72219.06
$ws.Range("K80").Value = 72219.06
$ws.Range("M80").Value = -71221.06
$ws.Range("H83").Value = 12866259
$ws.Range("I83").Value = 72219.06
$ws.Range("K83").Value = 361095.3
$ws.Range("M83").Value = -356103.3
$ws.Range("H122").Value = 5417.1113
$ws.Range("I122").Value = 3390
$ws.Range("K122").Value = 10170
$ws.Range("M122").Value = -7720
$ws.Range("H136").Value = 16752.316
$ws.Range("J136").Value = 16752.316
$ws.Range("L136").Value = 50256.948
$ws.Range("N136").Value = -55356.948

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4814.074
$ws.Range("I136").Value = 2748.8235
$ws.Range("K136").Value = 8246.470499999999
$ws.Range("M136").Value = -5696.470499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 45000
$ws.Range("J95").Value = 45000
$ws.Range("L95").Value = 45000
$ws.Range("N95").Value = -50492
$ws.Range("H126").Value = 1704
$ws.Range("I126").Value = 1461
$ws.Range("K126").Value = 4383
$ws.Range("M126").Value = -1913
